$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.7368846666666666
$ws.Range("H2").Value = 2.210654
$ws.Range("I2").Value = 0.1365512798651915
$ws.Range("J2").Value = 0.1426235012188444
$ws.Range("M2").Value = 12.36292333333333
$ws.Range("N2").Value = 37.08877
$ws.Range("O2").Value = 0.918273862214392
$ws.Range("P2").Value = 0.9303794466068031
$ws.Range("Q2").Value = 9.110048639508888
$ws.Range("R2").Value = 81.99043775557999
$ws.Range("S2").Value = 0.1253914711521277
$ws.Range("T2").Value = 0.1326939741371131
$ws.Range("G3").Value = 0.7368846666666666
$ws.Range("H3").Value = 2.210654
$ws.Range("I3").Value = 0.1365512798651915
$ws.Range("J3").Value = 0.1426235012188444
$ws.Range("O3").Value = 0.04269179184247177
$ws.Range("P3").Value = 0.04325459680761149
$ws.Range("Q3").Value = 0.4235384629753333
$ws.Range("R3").Value = 3.811846166778
$ws.Range("S3").Value = 0.00582961881582786
$ws.Range("T3").Value = 0.006169122040510999
$ws.Range("G4").Value = 0.7368846666666666
$ws.Range("H4").Value = 2.210654
$ws.Range("I4").Value = 0.1365512798651915
$ws.Range("J4").Value = 0.1426235012188444
$ws.Range("M4").Value = 0.525528
$ws.Range("N4").Value = 1.051056
$ws.Range("O4").Value = 0.03903434594313629
$ws.Range("P4").Value = 0.02636595658558534
$ws.Range("Q4").Value = 0.387253525104
$ws.Range("R4").Value = 2.323521150624
$ws.Range("S4").Value = 0.005330189897235904
$ws.Range("T4").Value = 0.003760405041220229
$ws.Range("I5").Value = 0.171453468119196
$ws.Range("J5").Value = 0.1790777350707694
$ws.Range("M5").Value = 12.36292333333333
$ws.Range("N5").Value = 37.08877
$ws.Range("O5").Value = 0.918273862214392
$ws.Range("P5").Value = 0.9303794466068031
$ws.Range("Q5").Value = 11.43855579764889
$ws.Range("R5").Value = 102.94700217884
$ws.Range("S5").Value = 0.1574412383598662
$ws.Range("T5").Value = 0.1666102440547421
$ws.Range("I6").Value = 0.171453468119196
$ws.Range("J6").Value = 0.1790777350707694
$ws.Range("O6").Value = 0.04269179184247177
$ws.Range("P6").Value = 0.04325459680761149
$ws.Range("S6").Value = 0.007319655771614583
$ws.Range("T6").Value = 0.007745935227706396
$ws.Range("I7").Value = 0.171453468119196
$ws.Range("J7").Value = 0.1790777350707694
$ws.Range("M7").Value = 0.525528
$ws.Range("N7").Value = 1.051056
$ws.Range("O7").Value = 0.03903434594313629
$ws.Range("P7").Value = 0.02636595658558534
$ws.Range("Q7").Value = 0.486234621792
$ws.Range("R7").Value = 2.917407730752
$ws.Range("S7").Value = 0.006692573987715184
$ws.Range("T7").Value = 0.004721555788320858
$ws.Range("G8").Value = 1.654700333333333
$ws.Range("H8").Value = 4.964101
$ws.Range("I8").Value = 0.3066306825627515
$ws.Range("J8").Value = 0.3202660683327045
$ws.Range("M8").Value = 12.36292333333333
$ws.Range("N8").Value = 37.08877
$ws.Range("O8").Value = 0.918273862214392
$ws.Range("P8").Value = 0.9303794466068031
$ws.Range("Q8").Value = 20.45693336064111
$ws.Range("R8").Value = 184.11240024577
$ws.Range("S8").Value = 0.281570941150333
$ws.Range("T8").Value = 0.2979689674223182
$ws.Range("G9").Value = 1.654700333333333
$ws.Range("H9").Value = 4.964101
$ws.Range("I9").Value = 0.3066306825627515
$ws.Range("J9").Value = 0.3202660683327045
$ws.Range("O9").Value = 0.04269179184247177
$ws.Range("P9").Value = 0.04325459680761149
$ws.Range("Q9").Value = 0.9510704558896667
$ws.Range("R9").Value = 8.559634103007001
$ws.Range("S9").Value = 0.01309061327248402
$ws.Range("T9").Value = 0.01385297965689009
$ws.Range("G10").Value = 1.654700333333333
$ws.Range("H10").Value = 4.964101
$ws.Range("I10").Value = 0.3066306825627515
$ws.Range("J10").Value = 0.3202660683327045
$ws.Range("M10").Value = 0.525528
$ws.Range("N10").Value = 1.051056
$ws.Range("O10").Value = 0.03903434594313629
$ws.Range("P10").Value = 0.02636595658558534
$ws.Range("Q10").Value = 0.8695913567760001
$ws.Range("R10").Value = 5.217548140656
$ws.Range("S10").Value = 0.01196912813993445
$ws.Range("T10").Value = 0.008444121253496195
$ws.Range("G11").Value = 0.6892575000000001
$ws.Range("H11").Value = 1.378515
$ws.Range("I11").Value = 0.1277255424616637
$ws.Range("J11").Value = 0.0889368647389846
$ws.Range("M11").Value = 12.36292333333333
$ws.Range("N11").Value = 37.08877
$ws.Range("O11").Value = 0.918273862214392
$ws.Range("P11").Value = 0.9303794466068031
$ws.Range("Q11").Value = 8.521237629425
$ws.Range("R11").Value = 51.12742577655
$ws.Range("S11").Value = 0.1172870271797002
$ws.Range("T11").Value = 0.08274503099880058
$ws.Range("G12").Value = 0.6892575000000001
$ws.Range("H12").Value = 1.378515
$ws.Range("I12").Value = 0.1277255424616637
$ws.Range("J12").Value = 0.0889368647389846
$ws.Range("O12").Value = 0.04269179184247177
$ws.Range("P12").Value = 0.04325459680761149
$ws.Range("Q12").Value = 0.3961638440175
$ws.Range("R12").Value = 2.376983064105
$ws.Range("S12").Value = 0.005452832271740134
$ws.Range("T12").Value = 0.003846928225617858
$ws.Range("G13").Value = 0.6892575000000001
$ws.Range("H13").Value = 1.378515
$ws.Range("I13").Value = 0.1277255424616637
$ws.Range("J13").Value = 0.0889368647389846
$ws.Range("M13").Value = 0.525528
$ws.Range("N13").Value = 1.051056
$ws.Range("O13").Value = 0.03903434594313629
$ws.Range("P13").Value = 0.02636595658558534
$ws.Range("Q13").Value = 0.3622241154600001
$ws.Range("R13").Value = 1.44889646184
$ws.Range("S13").Value = 0.004985683010223323
$ws.Range("T13").Value = 0.002344905514566144
$ws.Range("G14").Value = 1.390322
$ws.Range("H14").Value = 4.170966
$ws.Range("I14").Value = 0.2576390269911973
$ws.Range("J14").Value = 0.2690958306386971
$ws.Range("M14").Value = 12.36292333333333
$ws.Range("N14").Value = 37.08877
$ws.Range("O14").Value = 0.918273862214392
$ws.Range("P14").Value = 0.9303794466068031
$ws.Range("Q14").Value = 17.18844429464666
$ws.Range("R14").Value = 154.69599865182
$ws.Range("S14").Value = 0.2365831843723647
$ws.Range("T14").Value = 0.250361229993829
$ws.Range("G15").Value = 1.390322
$ws.Range("H15").Value = 4.170966
$ws.Range("I15").Value = 0.2576390269911973
$ws.Range("J15").Value = 0.2690958306386971
$ws.Range("O15").Value = 0.04269179184247177
$ws.Range("P15").Value = 0.04325459680761149
$ws.Range("Q15").Value = 0.799113985618
$ws.Range("R15").Value = 7.192025870562
$ws.Range("S15").Value = 0.01099907171080516
$ws.Range("T15").Value = 0.01163963165688615
$ws.Range("G16").Value = 1.390322
$ws.Range("H16").Value = 4.170966
$ws.Range("I16").Value = 0.2576390269911973
$ws.Range("J16").Value = 0.2690958306386971
$ws.Range("M16").Value = 0.525528
$ws.Range("N16").Value = 1.051056
$ws.Range("O16").Value = 0.03903434594313629
$ws.Range("P16").Value = 0.02636595658558534
$ws.Range("Q16").Value = 0.730653140016
$ws.Range("R16").Value = 4.383918840096
$ws.Range("S16").Value = 0.01005677090802742
$ws.Range("T16").Value = 0.007094968987981914
